$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1. Clear the old layout (rows 1-13, cols A-K) so we can rebuild it
#    with the new "cantrals by cantons" column layout.
# ---------------------------------------------------------------------
$ws.Range("A1:K13").Clear()

# ---------------------------------------------------------------------
# 2. Header row (row 1)
# ---------------------------------------------------------------------
$headers = @("idx", "idx2", "Name", "Date Start", "Date End", "(m3/s)", "(MW1)", "(MW2)", "(GWh) Winter", "(GWh) Summer", "(GWh) Year")
$ws.Range("A1:K1").Value = $headers

$ws.Range("A1:K1").Font.Name = "Arial"
$ws.Range("A1:K1").Font.Size = 9

# F1:K1 pick up the distinct "units" style present in the target file.
$ws.Range("F1:K1").Font.Name = "Arial"
$ws.Range("F1:K1").Font.Size = 9

# ---------------------------------------------------------------------
# 3. Data rows (rows 2-12) - one power plant per row
# ---------------------------------------------------------------------
$data = @(
    @(1, 207500, "Emmenhof", 1863, 1986, 12, 0.33, 0.32, 0.8, 0.9, 1.7),
    @(2, 207400, "Biberist (Papierfabrik)", 1864, 1985, 12, 0.5, 0.49, 1.5, 1.2, 2.7),
    @(3, 207600, "Untere Emmengasse", 1876, 2001, 13, 0.86, 0.82, 2.7, 2.5, 5.2),
    @(4, 207650, "Luterbach", 1888, 1988, 12, 0.32, 0.3, 0.77, 0.72, 1.49),
    @(5, 208400, "Aarau Stadt", 1893, 1964, 394, 13.88, 13.19, 39.28, 50.18, 89.46),
    @(6, 208300, "Gösgen", 1917, 2000, 380, 47.71, 45.57, 124.81, 156.98, 281.79),
    @(7, 208000, "Schwarzhäusern", 1923, 1979, 200, 0.8, 0.6, 1.8, 2.21, 4.02),
    @(8, 207700, "Flumenthal", 1970, 2009, 350, 14.53, 13.48, 40.37, 50.92, 91.29),
    @(9, 110450, "Dornachbrugg", 1996, $null, 20, 0.77, 0.77, 1.91, 1.56, 3.47),
    @(10, 207900, "Wynau", 1996, $null, 220, 1.18, 1.02, 2.16, 2.84, 5),
    @(11, 208100, "Ruppoldingen", 2000, $null, 475, 11.5, 10.75, 25.3, 32.2, 57.5)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    if ($null -ne $rec[4]) {
        $ws.Cells.Item($row, 5).Value = $rec[4]
    }
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $rec[7]
    $ws.Cells.Item($row, 9).Value = $rec[8]
    $ws.Cells.Item($row, 10).Value = $rec[9]
    $ws.Cells.Item($row, 11).Value = $rec[10]
    $row = $row + 1
}

# Fonts / number formats for the data block
$ws.Range("A2:K12").Font.Name = "Arial"
$ws.Range("A2:K12").Font.Size = 9

$ws.Range("A2:B12").NumberFormat = "0"
$ws.Range("D2:E12").NumberFormat = "0"
$ws.Range("F2:K12").NumberFormat = "0.00"

# ---------------------------------------------------------------------
# 4. Sheet-level bookkeeping to match the edited workbook
# ---------------------------------------------------------------------
$ws.Range("A2:K2").Select()
